$wb = $excel.ActiveWorkbook

# --- Sheet "Daily" ---
$daily = $wb.Worksheets.Item("Daily")
$daily.Range("G2").Value = 3243.21
$daily.Range("H2").Value = 6521.6
$daily.Range("I2").Value = 772.8099999999999
$daily.Range("J2").Value = 3242.83
$daily.Range("K2").Value = 6198.49
$daily.Range("L2").Value = 777.37

# --- Sheet "Hourly" ---
$hourly = $wb.Worksheets.Item("Hourly")

# Row 9
$hourly.Range("I9").Value = 69.28
$hourly.Range("L9").Value = 26.84

# Row 10
$hourly.Range("H10").Value = 119.58
$hourly.Range("I10").Value = 454.9
$hourly.Range("K10").Value = 119.58
$hourly.Range("L10").Value = 420.16

# Row 11
$hourly.Range("H11").Value = 268.31
$hourly.Range("I11").Value = 651.85
$hourly.Range("J11").Value = 74.72
$hourly.Range("K11").Value = 268.31
$hourly.Range("L11").Value = 632.89

# Row 12
$hourly.Range("H12").Value = 395.86
$hourly.Range("I12").Value = 747.33
$hourly.Range("K12").Value = 395.86
$hourly.Range("L12").Value = 730.5

# Row 13
$hourly.Range("I13").Value = 795.59
$hourly.Range("J13").Value = 96.06999999999999

# Row 14
$hourly.Range("J14").Value = 98.98999999999999

# Row 15
$hourly.Range("H15").Value = 500.98
$hourly.Range("J15").Value = 97.53
$hourly.Range("K15").Value = 500.98

# Row 16
$hourly.Range("M16").Value = 93.23

# Row 17
$hourly.Range("I17").Value = 690.98

# Row 18
$hourly.Range("I18").Value = 539.42
$hourly.Range("K18").Value = 169.49
$hourly.Range("L18").Value = 515.48
$hourly.Range("M18").Value = 56.54

# Row 19
$hourly.Range("I19").Value = 187.08
$hourly.Range("K19").Value = 32.6
$hourly.Range("L19").Value = 120.12
$hourly.Range("M19").Value = 20.05
